$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reorder the "Periodo Mora" column (E16:E18) from descending (2404,2403,2402)
# to ascending (2402,2403,2404).
$ws.Range("E16").Value = "2402"
$ws.Range("E17").Value = "2403"
$ws.Range("E18").Value = "2404"

# Update "Salario Basico" column (G16:G18) with the new value.
$ws.Range("G16").Value = 1423500
$ws.Range("G17").Value = 1423500
$ws.Range("G18").Value = 1423500
